# Reorders / updates rows 2-10 of "CATEGORIA D - SIN VENTAS" so that the
# article list reflects the refreshed ABC+D classification export: several
# articles swapped rank position and picked up new rotation/discount figures
# (and, for the now-unique-SKU rows, a fresh Talla/Color pair) while others
# simply got refreshed stock/rotation numbers and purchase-origin notes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CATEGORIA D – SIN VENTAS")

# Columns A (Articulo), E (Familia) and W (Escenario) hold digit-only codes
# that must stay text (e.g. "3201020001", "32", "14") instead of being
# auto-coerced to numbers when assigned through .Value.
$ws.Range("A2:A10").NumberFormat = "@"
$ws.Range("E2:E10").NumberFormat = "@"
$ws.Range("W2:W10").NumberFormat = "@"

# Each inner array is:
#   row, A Articulo, B Nombre articulo, C Talla, D Color, E Familia, F Nombre Familia,
#   G Rotacion Familia, H Ventas, I Importe ventas, J Beneficio, K Tasa de venta,
#   L Rotacion excedida, M Stock minimo, N Stock maximo, O Stock Final,
#   P Antiguedad Ultima Venta, Q Antiguedad Stock, R % Rotacion Consumido,
#   S Descuento Sugerido, T Riesgo, U Accion Sugerida, V Origen Stock Final, W Escenario
$rows = @(
    @(2, "3201020001", "BIG BAG MANTILLO 500L (NO VENTA)", "", "", "32", "MANTENIMIENTO", 90, 0, 0, 0, 0, 7, 0, 0, 7, 92, 12, 13.33, 0, "Crítico", "LIQUIDACIÓN URGENTE: Aplicar descuento 0% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima.", "Compra 19/05/2025", "14"),
    @(3, "3203050008", "BIG BAG MARMOLINA BLANCA 500KG (NO VENTA)", "12I18", "UNICO", "32", "MANTENIMIENTO", 90, 0, 0, 0, 0, 3, 0, 0, 3, 92, 12, 13.33, 0, "Crítico", "LIQUIDACIÓN URGENTE: Aplicar descuento 0% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima.", "Compra 19/05/2025", "14"),
    @(4, "3102110006", "BIG BAG TIERRA ENRIQUECIDA 650L (NO VENTA)", "", "", "31", "TIERRAS", 90, 0, 0, 0, 0, 10, 0, 0, 10, 92, 92, 102.22, 20, "Crítico", "LIQUIDACIÓN URGENTE: Aplicar descuento 20% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 144.76€. Prioridad máxima.", "Stock inicial", "14"),
    @(5, "3203050005", "BIG BAG GRAVA VOLCANICA MARRON 500L (NO VENTA)", "5I10", "UNICO", "32", "MANTENIMIENTO", 90, 0, 0, 0, 0, 3, 0, 0, 3, 92, 12, 13.33, 0, "Crítico", "LIQUIDACIÓN URGENTE: Aplicar descuento 0% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima.", "Compra 19/05/2025", "14"),
    @(6, "3102110006", "BIG BAG TIERRA ENRIQUECIDA 500L (NO VENTA)", "", "", "31", "TIERRAS", 90, 0, 0, 0, 0, 7, 0, 0, 7, 92, 12, 13.33, 0, "Crítico", "LIQUIDACIÓN URGENTE: Aplicar descuento 0% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima.", "Compra 19/05/2025", "14"),
    @(7, "3202030010", "BIG BAG RECEBO CESPED 650L (NO VENTA)", "", "", "32", "MANTENIMIENTO", 90, 0, 0, 0, 0, 1, 0, 0, 1, 92, 92, 102.22, 20, "Crítico", "LIQUIDACIÓN URGENTE: Aplicar descuento 20% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 16.0€. Prioridad máxima.", "Stock inicial", "14"),
    @(8, "3102110001", "BIG BAG JABRE", "1M3", "UNICO", "31", "TIERRAS", 90, 0, 0, 0, 0, 14, 0, 0, 14, 92, 92, 102.22, 20, "Crítico", "LIQUIDACIÓN URGENTE: Aplicar descuento 20% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 458.74€. Prioridad máxima.", "Stock inicial", "14"),
    @(9, "3101010010", "COMPO BIO SUBSTRATO AQUA DEPOT", "", "", "31", "TIERRAS", 90, 0, 0, 0, 0, 41, 0, 0, 41, 92, 92, 102.22, 20, "Crítico", "LIQUIDACIÓN URGENTE: Aplicar descuento 20% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 244.24€. Prioridad máxima.", "Stock inicial", "14"),
    @(10, "3203050002", "BIG BAG CANTO RODADO BLANCO 500KG", "20I40", "UNICO", "32", "MANTENIMIENTO", 90, 0, 0, 0, 0, 2, 0, 0, 2, 92, 50, 55.56, 0, "Crítico", "LIQUIDACIÓN URGENTE: Aplicar descuento 0% inmediato. Eliminar del catálogo próxima temporada. Capital liberado: 0.0€. Prioridad máxima.", "Compra 11/04/2025", "14")
)

foreach ($r in $rows) {
    $row = $r[0]
    $ws.Cells.Item($row, 1).Value  = $r[1]   # A Articulo
    $ws.Cells.Item($row, 2).Value  = $r[2]   # B Nombre articulo
    $ws.Cells.Item($row, 3).Value  = $r[3]   # C Talla
    $ws.Cells.Item($row, 4).Value  = $r[4]   # D Color
    $ws.Cells.Item($row, 5).Value  = $r[5]   # E Familia
    $ws.Cells.Item($row, 6).Value  = $r[6]   # F Nombre Familia
    $ws.Cells.Item($row, 7).Value  = $r[7]   # G Rotacion Familia (dias)
    $ws.Cells.Item($row, 8).Value  = $r[8]   # H Ventas (unidades)
    $ws.Cells.Item($row, 9).Value  = $r[9]   # I Importe ventas (EUR)
    $ws.Cells.Item($row, 10).Value = $r[10]  # J Beneficio (importe EUR)
    $ws.Cells.Item($row, 11).Value = $r[11]  # K Tasa de venta (%)
    $ws.Cells.Item($row, 12).Value = $r[12]  # L Rotacion excedida (unidades)
    $ws.Cells.Item($row, 13).Value = $r[13]  # M Stock minimo (unidades)
    $ws.Cells.Item($row, 14).Value = $r[14]  # N Stock maximo (unidades)
    $ws.Cells.Item($row, 15).Value = $r[15]  # O Stock Final (unidades)
    $ws.Cells.Item($row, 16).Value = $r[16]  # P Antiguedad Ultima Venta (dias)
    $ws.Cells.Item($row, 17).Value = $r[17]  # Q Antiguedad Stock (dias)
    $ws.Cells.Item($row, 18).Value = $r[18]  # R % Rotacion Consumido
    $ws.Cells.Item($row, 19).Value = $r[19]  # S Descuento Sugerido (%)
    $ws.Cells.Item($row, 20).Value = $r[20]  # T Riesgo de Merma / inmovilizado
    $ws.Cells.Item($row, 21).Value = $r[21]  # U Accion Sugerida
    $ws.Cells.Item($row, 22).Value = $r[22]  # V Origen Stock Final
    $ws.Cells.Item($row, 23).Value = $r[23]  # W Escenario
}
